$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.547.20'
$ws.Range("E2").Value = '  +0.16%  '

$ws.Range("D3").Value = '2.641.16'
$ws.Range("E3").Value = '  +1.19%  '

$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").Value = '''536.37'
$ws.Range("E5").Value = '  -0.22%  '

$ws.Range("D6").Value = '''145.23'
$ws.Range("E6").Value = '  +2.97%  '

$ws.Range("D7").Value = '''0.998'
$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("D8").Value = '''0.572'
$ws.Range("E8").Value = '  +0.44%  '

$ws.Range("D9").Value = '2.653.46'
$ws.Range("E9").Value = '  +1.33%  '

$ws.Range("D10").Value = '''6.63'
$ws.Range("E10").Value = '  +2.46%  '

$ws.Range("D11").Value = '''0.103'
$ws.Range("E11").Value = '  -0.79%  '

$ws.Range("D12").Value = '''0.338'
$ws.Range("E12").Value = '  -0.01%  '

$ws.Range("E13").Value = '  -0.34%  '

$ws.Range("D14").Value = '3.108.79'
$ws.Range("E14").Value = '  +1.24%  '

$ws.Range("D15").Value = '59.422.73'
$ws.Range("E15").Value = '  +0.08%  '

$ws.Range("D16").Value = '''21.12'
$ws.Range("E16").Value = '  +2.42%  '

$ws.Range("D17").Value = '2.647.74'
$ws.Range("E17").Value = '  -0.68%  '

$ws.Range("D18").Value = '''0.0000135'
$ws.Range("E18").Value = '  +0.49%  '

$ws.Range("D19").Value = '''340.41'
$ws.Range("E19").Value = '  -1.37%  '

$ws.Range("E20").Value = '  +0.83%  '

$ws.Range("D21").Value = '''10.37'
$ws.Range("E21").Value = '  +2.10%  '

$ws.Range("D22").Value = '''6.30'
$ws.Range("E22").Value = '  -1.54%  '

$ws.Range("D23").Value = '''1.00'
$ws.Range("E23").Value = '  +0.08%  '

$ws.Range("D24").Value = '''67.03'
$ws.Range("E24").Value = '  -0.13%  '

$ws.Range("E25").Value = '  +1.44%  '

$ws.Range("E26").Value = '  -1.53%  '

$ws.Range("D27").Value = '''0.999'
$ws.Range("E27").Value = '  +0.04%  '

$ws.Range("D28").Value = '''7.28'
$ws.Range("E28").Value = '  +0.62%  '

$ws.Range("D29").Value = '0.0₃0747'
$ws.Range("E29").Value = '  -0.37%  '

$ws.Range("D30").Value = '''0.998'
$ws.Range("E30").Value = '  -0.05%  '

$ws.Range("E31").Value = '  +0.34%  '

$ws.Range("D32").Value = '''5.84'
$ws.Range("E32").Value = '  -0.48%  '

$ws.Range("D33").Value = '''18.89'
$ws.Range("E33").Value = '  -0.08%  '

$ws.Range("D34").Value = '''150.79'
$ws.Range("E34").Value = '  +1.23%  '

$ws.Range("D35").Value = '''4.00'
$ws.Range("E35").Value = '  -0.28%  '

$ws.Range("E36").Value = '  +0.91%  '

$ws.Range("D37").Value = '''0.837'
$ws.Range("E37").Value = '  -0.65%  '

$ws.Range("D38").Value = '''0.835'
$ws.Range("E38").Value = '  -0.66%  '

$ws.Range("E39").Value = '  -0.91%  '

$ws.Range("D40").Value = '''289.34'
$ws.Range("E40").Value = '  +4.35%  '

$ws.Range("D41").Value = '''3.59'
$ws.Range("E41").Value = '  +0.80%  '

$ws.Range("D42").Value = '''0.998'
$ws.Range("E42").Value = '  -0.05%  '

$ws.Range("E43").Value = '  +0.77%  '

$ws.Range("E44").Value = '  -0.09%  '

$ws.Range("D45").Value = '''19.34'
$ws.Range("E45").Value = '  +3.22%  '

$ws.Range("E46").Value = '  +1.57%  '

$ws.Range("D47").Value = '''0.0946'
$ws.Range("E47").Value = '  -1.75%  '

$ws.Range("D48").Value = '1.973.00'
$ws.Range("E48").Value = '  +1.00%  '

$ws.Range("D49").Value = '''0.0226'
$ws.Range("E49").Value = '  +0.80%  '

$ws.Range("E50").Value = '  +0.54%  '

$ws.Range("D51").Value = '''18.31'
$ws.Range("E51").Value = '  -0.47%  '
